$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 20.94835401735078
$ws.Range("E2").Value = 19.77249717712402
$ws.Range("F2").Value = 21.04785020659737
$ws.Range("G2").Value = 19.75440664722321
$ws.Range("H2").Value = 174270076
$ws.Range("I2").Value = "G"

$ws.Range("D3").Value = 19.36546554705393
$ws.Range("E3").Value = 20.08906936645508
$ws.Range("F3").Value = 20.24283539372877
$ws.Range("G3").Value = 19.17552019141397
$ws.Range("H3").Value = 174270076
$ws.Range("I3").Value = "G"

$ws.Range("D4").Value = 21.30110907043843
$ws.Range("E4").Value = 22.41365242004395
$ws.Range("F4").Value = 22.70309397515155
$ws.Range("G4").Value = 20.87599200202511
$ws.Range("H4").Value = 174270076
$ws.Range("I4").Value = "G"

$ws.Range("D5").Value = 22.32320303014261
$ws.Range("E5").Value = 21.63577842712402
$ws.Range("F5").Value = 22.32320303014261
$ws.Range("G5").Value = 20.66795751406097
$ws.Range("H5").Value = 174270076
$ws.Range("I5").Value = "G"

$ws.Range("D6").Value = 24.50306156736417
$ws.Range("E6").Value = 25.22666549682617
$ws.Range("F6").Value = 25.79650510083437
$ws.Range("G6").Value = 24.31311445762548
$ws.Range("H6").Value = 174270076
$ws.Range("I6").Value = "G"

$ws.Range("D7").Value = 24.26788527062729
$ws.Range("E7").Value = 24.21361541748047
$ws.Range("F7").Value = 24.83772304168768
$ws.Range("G7").Value = 23.28197402488788
$ws.Range("H7").Value = 174270076
$ws.Range("I7").Value = "G"

$ws.Range("D8").Value = 21.71718242502667
$ws.Range("E8").Value = 20.79458618164062
$ws.Range("F8").Value = 21.78049754904591
$ws.Range("G8").Value = 20.74031632179667
$ws.Range("H8").Value = 174270076
$ws.Range("I8").Value = "G"

$ws.Range("D9").Value = 22.16039082253095
$ws.Range("E9").Value = 22.32320213317871
$ws.Range("F9").Value = 22.51314751212923
$ws.Range("G9").Value = 21.39156059180974
$ws.Range("H9").Value = 174270076
$ws.Range("I9").Value = "G"

$ws.Range("D10").Value = 22.43362722353422
$ws.Range("E10").Value = 22.14345932006836
$ws.Range("F10").Value = 22.53337308892703
$ws.Range("G10").Value = 21.16414178110251
$ws.Range("H10").Value = 174270076
$ws.Range("I10").Value = "G"

$ws.Range("D11").Value = 25.30901777470153
$ws.Range("E11").Value = 26.35409355163575
$ws.Range("F11").Value = 26.46314573589201
$ws.Range("G11").Value = 24.83646119846714
$ws.Range("H11").Value = 174270076
$ws.Range("I11").Value = "G"

$ws.Range("D12").Value = 26.18290791734447
$ws.Range("E12").Value = 27.73111534118652
$ws.Range("F12").Value = 27.82218493561392
$ws.Range("G12").Value = 26.18290791734447
$ws.Range("H12").Value = 174270076
$ws.Range("I12").Value = "G"

$ws.Range("D13").Value = 29.1427989848308
$ws.Range("E13").Value = 30.96764373779297
$ws.Range("F13").Value = 31.74320558580439
$ws.Range("G13").Value = 29.04243187080192
$ws.Range("H13").Value = 174270076
$ws.Range("I13").Value = "G"

$ws.Range("D14").Value = 29.2673206876934
$ws.Range("E14").Value = 29.16671371459961
$ws.Range("F14").Value = 29.91668852810683
$ws.Range("G14").Value = 28.50819962713508
$ws.Range("H14").Value = 174270076
$ws.Range("I14").Value = "G"

$ws.Range("D15").Value = 26.43185200831845
$ws.Range("E15").Value = 27.85291862487793
$ws.Range("F15").Value = 28.32049582257778
$ws.Range("G15").Value = 26.31266648672007
$ws.Range("H15").Value = 174270076
$ws.Range("I15").Value = "G"

$ws.Range("D16").Value = 28.19744042515341
$ws.Range("E16").Value = 25.19204139709473
$ws.Range("F16").Value = 28.26177589467594
$ws.Range("G16").Value = 24.31891266012022
$ws.Range("H16").Value = 174270076
$ws.Range("I16").Value = "G"

$ws.Range("D17").Value = 24.57625605043068
$ws.Range("E17").Value = 27.48817825317383
$ws.Range("F17").Value = 27.58032804594137
$ws.Range("G17").Value = 24.17079801682017
$ws.Range("H17").Value = 174270076
$ws.Range("I17").Value = "G"

$ws.Range("D18").Value = 32.65801148820098
$ws.Range("E18").Value = 33.53566741943359
$ws.Range("F18").Value = 33.74815194607746
$ws.Range("G18").Value = 32.38085791102198
$ws.Range("H18").Value = 174270076
$ws.Range("I18").Value = "G"

$ws.Range("D19").Value = 35.49306092889683
$ws.Range("E19").Value = 36.74314117431641
$ws.Range("F19").Value = 37.79876512064545
$ws.Range("G19").Value = 35.14118510266787
$ws.Range("H19").Value = 174270076
$ws.Range("I19").Value = "G"

$ws.Range("D20").Value = 36.08675677915155
$ws.Range("E20").Value = 36.34657287597656
$ws.Range("F20").Value = 36.51359852100088
$ws.Range("G20").Value = 34.66704233617366
$ws.Range("H20").Value = 174270076
$ws.Range("I20").Value = "G"

$ws.Range("D21").Value = 39.47280558900383
$ws.Range("E21").Value = 41.1651611328125
$ws.Range("F21").Value = 41.90905309484543
$ws.Range("G21").Value = 38.99857202453762
$ws.Range("H21").Value = 174270076
$ws.Range("I21").Value = "G"

$ws.Range("D22").Value = 26.12987068002252
$ws.Range("E22").Value = 32.09601974487305
$ws.Range("F22").Value = 33.58755567754736
$ws.Range("G22").Value = 24.17222754582363
$ws.Range("H22").Value = 174270076
$ws.Range("I22").Value = "G"

$ws.Range("D23").Value = 34.10683361410414
$ws.Range("E23").Value = 37.2193489074707
$ws.Range("F23").Value = 37.5651831477206
$ws.Range("G23").Value = 32.26549350000506
$ws.Range("H23").Value = 174270076
$ws.Range("I23").Value = "G"

$ws.Range("D24").Value = 36.5919817057332
$ws.Range("E24").Value = 32.20656585693359
$ws.Range("F24").Value = 37.68833745522294
$ws.Range("G24").Value = 31.77552159540893
$ws.Range("H24").Value = 174270076
$ws.Range("I24").Value = "G"

$ws.Range("D25").Value = 38.82200684566785
$ws.Range("E25").Value = 35.95708465576172
$ws.Range("F25").Value = 39.84585981969669
$ws.Range("G25").Value = 35.95708465576172
$ws.Range("H25").Value = 174270076
$ws.Range("I25").Value = "G"

$ws.Range("D26").Value = 40.64685389713144
$ws.Range("E26").Value = 44.76239395141602
$ws.Range("F26").Value = 46.14679907323642
$ws.Range("G26").Value = 40.41141104471697
$ws.Range("H26").Value = 174270076
$ws.Range("I26").Value = "G"

$ws.Range("D27").Value = 42.96178698157817
$ws.Range("E27").Value = 47.02102279663086
$ws.Range("F27").Value = 47.36086252898332
$ws.Range("G27").Value = 42.8957067335288
$ws.Range("H27").Value = 174270076
$ws.Range("I27").Value = "G"

$ws.Range("D28").Value = 45.31251845885882
$ws.Range("E28").Value = 46.68418884277344
$ws.Range("F28").Value = 48.71805065330654
$ws.Range("G28").Value = 44.3665368927752
$ws.Range("H28").Value = 174270076
$ws.Range("I28").Value = "G"

$ws.Range("D29").Value = 50.58405119667986
$ws.Range("E29").Value = 47.1618537902832
$ws.Range("F29").Value = 51.21919402449004
$ws.Range("G29").Value = 44.31792290845708
$ws.Range("H29").Value = 174270076
$ws.Range("I29").Value = "G"

$ws.Range("D30").Value = 41.61108006623338
$ws.Range("E30").Value = 38.29246520996094
$ws.Range("F30").Value = 42.1150558710759
$ws.Range("G30").Value = 38.23541031864784
$ws.Range("H30").Value = 174270076
$ws.Range("I30").Value = "G"

$ws.Range("D31").Value = 40.45002480611166
$ws.Range("E31").Value = 45.8471794128418
$ws.Range("F31").Value = 46.03788855632935
$ws.Range("G31").Value = 39.60135674994331
$ws.Range("H31").Value = 174270076
$ws.Range("I31").Value = "G"

$ws.Range("D32").Value = 42.2608406204096
$ws.Range("E32").Value = 46.37218856811523
$ws.Range("F32").Value = 46.61122046795088
$ws.Range("G32").Value = 40.21472814139045
$ws.Range("H32").Value = 174270076
$ws.Range("I32").Value = "G"

$ws.Range("D33").Value = 44.76733024355548
$ws.Range("E33").Value = 45.33303451538086
$ws.Range("F33").Value = 46.47403331281059
$ws.Range("G33").Value = 43.51127766165907
$ws.Range("H33").Value = 174270076
$ws.Range("I33").Value = "G"

$ws.Range("D34").Value = 44.32269747854671
$ws.Range("E34").Value = 42.84174728393555
$ws.Range("F34").Value = 44.41886155105944
$ws.Range("G34").Value = 41.44734622195554
$ws.Range("H34").Value = 174270076
$ws.Range("I34").Value = "G"

$ws.Range("D35").Value = 36.28250799732277
$ws.Range("E35").Value = 34.83468246459961
$ws.Range("F35").Value = 38.20329149082021
$ws.Range("G35").Value = 34.79607290004441
$ws.Range("H35").Value = 174270076
$ws.Range("I35").Value = "G"

$ws.Range("D36").Value = 35.0218216313244
$ws.Range("E36").Value = 32.4932746887207
$ws.Range("F36").Value = 35.9906133368034
$ws.Range("G36").Value = 31.59229810697348
$ws.Range("H36").Value = 174270076
$ws.Range("I36").Value = "G"

$ws.Range("D37").Value = 33.66524865844922
$ws.Range("E37").Value = 34.92003631591797
$ws.Range("F37").Value = 35.89273923602079
$ws.Range("G37").Value = 32.85790345369164
$ws.Range("H37").Value = 174270076
$ws.Range("I37").Value = "G"

$ws.Range("D38").Value = 32.242142750465
$ws.Range("E38").Value = 30.03404426574707
$ws.Range("F38").Value = 32.25191145683284
$ws.Range("G38").Value = 29.53575638143259
$ws.Range("H38").Value = 174270076
$ws.Range("I38").Value = "G"

$ws.Range("D39").Value = 31.59574044072415
$ws.Range("E39").Value = 34.02995681762695
$ws.Range("F39").Value = 34.52072633904141
$ws.Range("G39").Value = 29.81915529740221
$ws.Range("H39").Value = 174270076
$ws.Range("I39").Value = "G"

$ws.Range("D40").Value = 38.63818800200151
$ws.Range("E40").Value = 37.61335372924805
$ws.Range("F40").Value = 39.30827166041212
$ws.Range("G40").Value = 37.26845938700718
$ws.Range("H40").Value = 174270076
$ws.Range("I40").Value = "G"

$ws.Range("D41").Value = 42.69494131104493
$ws.Range("E41").Value = 48.14304351806641
$ws.Range("F41").Value = 48.19248486464085
$ws.Range("G41").Value = 41.81493965659031
$ws.Range("H41").Value = 174270076
$ws.Range("I41").Value = "G"

$ws.Range("D42").Value = 49.76499341349398
$ws.Range("E42").Value = 49.86420440673828
$ws.Range("F42").Value = 50.31066144565122
$ws.Range("G42").Value = 44.01066875821829
$ws.Range("H42").Value = 174270076
$ws.Range("I42").Value = "G"

$ws.Range("D43").Value = 43.88580945541839
$ws.Range("E43").Value = 43.87584686279297
$ws.Range("F43").Value = 46.56517312967035
$ws.Range("G43").Value = 43.5272321132193
$ws.Range("H43").Value = 174270076
$ws.Range("I43").Value = "G"
